$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, pushing existing rows 10-16 down to 11-17.
$ws.Rows.Item(10).Insert()

# Populate the new row 10 with this week's data (same market/category as the
# rows around it, new date and updated prices/origin).
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Terminal La Palmera de La Serena"
$ws.Range("C10").Value = "Coquimbo"
$ws.Range("D10").Value = 44839
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 100112013
$ws.Range("G10").Value = "Alcachofa"
$ws.Range("H10").Value = "Española"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 400
$ws.Range("K10").Value = 12000
$ws.Range("L10").Value = 13000
$ws.Range("M10").Value = 12500
$ws.Range("N10").Value = "`$/caja 30 unidades"
$ws.Range("O10").Value = "Provincia del Elquí"
$ws.Range("P10").Value = 417
$ws.Range("Q10").Value = 30
$ws.Range("R10").Value = "Hortaliza"

# Ensure the date cell keeps the same date-time number format used elsewhere
# in column D.
$ws.Range("D10").NumberFormat = $ws.Range("D11").NumberFormat
